$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vector Entrada")
$ws.Rows.Item(367).Delete()
$ws.Rows.Item(336).Delete()
